$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "270.76"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "3.63%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "26.73"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-1.39%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.721"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.03%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06127"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-1.26%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.750"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.51%"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8554"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.41%"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8993"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-1.22%"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1428"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "1.53%"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.05088"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "5.38%"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07161"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "0.96%"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03157"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-0.37%"

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.16%"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001529"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.79%"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0006093"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-1.01%"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006078"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.92%"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.465"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.04%"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.182"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.30%"

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "4.60%"

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "0.04%"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.841"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-6.51%"

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.30%"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001177"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-3.57%"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004152"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "0.68%"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001202"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.01%"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001678"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "3.81%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03968"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "1.44%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1121"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "0.75%"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.004198"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "1.56%"

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-6.63%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01170"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-13.08%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005160"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-0.28%"

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.00%"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.9061"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "433.70%"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.02994"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-16.63%"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002102"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.00%"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0002002"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.00%"
